$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 710-711, pushing the existing rows (710..795)
# down to (712..797). Excel's default insert behaviour copies formatting
# from the row above, so the date style on column D is preserved.
$ws.Range("A710:A711").EntireRow.Insert()

# Populate the new row 710 (Camote, "1a nueva(o)")
$ws.Cells.Item(710, 1).Value = 3
$ws.Cells.Item(710, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(710, 3).Value = "Coquimbo"
$ws.Cells.Item(710, 4).Value = 44946
$ws.Cells.Item(710, 5).Value = 5
$ws.Cells.Item(710, 6).Value = 100112045
$ws.Cells.Item(710, 7).Value = "Zapallo"
$ws.Cells.Item(710, 8).Value = "Camote"
$ws.Cells.Item(710, 9).Value = "1a nueva(o)"
$ws.Cells.Item(710, 10).Value = 120
$ws.Cells.Item(710, 11).Value = 700
$ws.Cells.Item(710, 12).Value = 700
$ws.Cells.Item(710, 13).Value = 700
$ws.Cells.Item(710, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(710, 15).Value = "Provincia de Talca"
$ws.Cells.Item(710, 16).Value = 700
$ws.Cells.Item(710, 17).Value = 1
$ws.Cells.Item(710, 18).Value = "Hortaliza"

# Populate the new row 711 (Paine, "1a nueva(o)")
$ws.Cells.Item(711, 1).Value = 3
$ws.Cells.Item(711, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(711, 3).Value = "Coquimbo"
$ws.Cells.Item(711, 4).Value = 44946
$ws.Cells.Item(711, 5).Value = 5
$ws.Cells.Item(711, 6).Value = 100112045
$ws.Cells.Item(711, 7).Value = "Zapallo"
$ws.Cells.Item(711, 8).Value = "Paine"
$ws.Cells.Item(711, 9).Value = "1a nueva(o)"
$ws.Cells.Item(711, 10).Value = 220
$ws.Cells.Item(711, 11).Value = 400
$ws.Cells.Item(711, 12).Value = 450
$ws.Cells.Item(711, 13).Value = 418
$ws.Cells.Item(711, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(711, 15).Value = "Provincia de Talca"
$ws.Cells.Item(711, 16).Value = 418
$ws.Cells.Item(711, 17).Value = 1
$ws.Cells.Item(711, 18).Value = "Hortaliza"
